$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 (the standalone value-only row holding the professor's name),
# shifting everything below it up by one row.
$ws.Rows.Item(13).Delete()

# Update cell contents that were overwritten with new values after the shift.
$ws.Range("B10").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C10").Value = "101761 - Arnaldo Márcio Ramalho Prata"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

$ws.Range("B18").Value = "101761 - Arnaldo Márcio Ramalho Prata"
$ws.Range("C18").Value = "101761 - Arnaldo Márcio Ramalho Prata"

$ws.Range("B19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."
$ws.Range("C19").Value = "Supervisão das atividades desenvolvidas pelo aluno durante o estágio."

$ws.Range("B20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio"
$ws.Range("C20").Value = "MF = Nota baseada em relatório final e no desempenho no estágio, a ser atribuída pelo professor orientador do estágio"

$ws.Range("B21").Value = "Não será oferecida recuperação."
$ws.Range("C21").Value = "Não será oferecida recuperação."
